# Rebuild rows 2-31 of Sheet1 (Wnt11-Fzd8 LR-pairs) with refreshed TPM-derived values.
# Columns A:T = Sending cluster, Ligand symbol, Receptor symbol, Target cluster,
# Ligand-expressing cells, Ligand detection rate, Ligand avg/total expression value,
# Ligand derived specificity (avg/total), Receptor-expressing cells, Receptor detection rate,
# Receptor avg/total expression value, Receptor derived specificity (avg/total),
# Edge avg/total expression weight, Edge avg/total expression derived specificity.
$data = @(
  @("ECs", "Wnt11", "Fzd8", "ECs", 1, 0.3333333333333333, 0.09047666666666666, 0.27143, 0.01063327555148536, 0.01063327555148536, 3, 1, 2.157506, 6.472517999999999, 0.3549648016839517, 0.3549648016839516, 0.1952039511933333, 1.75683556074, 0.003774438547383812, 0.003774438547383812),
  @("ECs", "Wnt11", "Fzd8", "FAPs", 1, 0.3333333333333333, 0.09047666666666666, 0.27143, 0.01063327555148536, 0.01063327555148536, 3, 1, 2.913654666666667, 8.740964, 0.4793705560628122, 0.4793705560628121, 0.2636177620577778, 2.37255985852, 0.005097279213884641, 0.005097279213884642),
  @("ECs", "Wnt11", "Fzd8", "Inflammatory-Mac", 1, 0.3333333333333333, 0.09047666666666666, 0.27143, 0.01063327555148536, 0.01063327555148536, 2, 0.6666666666666666, 0.018986, 0.05695799999999999, 0.003123681567871193, 0.003123681567871192, 0.001717789993333333, 0.01546010994, 0.0000332149668462702, 0.00003321496684627021),
  @("ECs", "Wnt11", "Fzd8", "MuSCs", 1, 0.3333333333333333, 0.09047666666666666, 0.27143, 0.01063327555148536, 0.01063327555148536, 3, 1, 0.9848966666666668, 2.95469, 0.1620406385718132, 0.1620406385718132, 0.08911016741111112, 0.8019915067000001, 0.001723022760472737, 0.001723022760472737),
  @("ECs", "Wnt11", "Fzd8", "Neutrophils", 1, 0.3333333333333333, 0.09047666666666666, 0.27143, 0.01063327555148536, 0.01063327555148536, 1, 0.3333333333333333, 0.003041, 0.009123, 0.0005003221135518961, 0.000500322113551896, 0.0002751395433333333, 0.00247625589, 0.000005320062897898857, 0.000005320062897898857),
  @("FAPs", "Wnt11", "Fzd8", "ECs", 3, 1, 4.990556000000001, 14.971668, 0.5865153863219087, 0.5865153863219088, 3, 1, 2.157506, 6.472517999999999, 0.3549648016839517, 0.3549648016839516, 10.767154513336, 96.90439062002399, 0.2081923177903426, 0.2081923177903426),
  @("FAPs", "Wnt11", "Fzd8", "FAPs", 3, 1, 4.990556000000001, 14.971668, 0.5865153863219087, 0.5865153863219088, 3, 1, 2.913654666666667, 8.740964, 0.4793705560628122, 0.4793705560628121, 14.54075677866133, 130.866811007952, 0.2811582068805285, 0.2811582068805285),
  @("FAPs", "Wnt11", "Fzd8", "Inflammatory-Mac", 3, 1, 4.990556000000001, 14.971668, 0.5865153863219087, 0.5865153863219088, 2, 0.6666666666666666, 0.018986, 0.05695799999999999, 0.003123681567871193, 0.003123681567871192, 0.094750696216, 0.852756265944, 0.001832087301526598, 0.001832087301526598),
  @("FAPs", "Wnt11", "Fzd8", "MuSCs", 3, 1, 4.990556000000001, 14.971668, 0.5865153863219087, 0.5865153863219088, 3, 1, 0.9848966666666668, 2.95469, 0.1620406385718132, 0.1620406385718132, 4.915181969213334, 44.23663772292001, 0.09503932773179581, 0.09503932773179581),
  @("FAPs", "Wnt11", "Fzd8", "Neutrophils", 3, 1, 4.990556000000001, 14.971668, 0.5865153863219087, 0.5865153863219088, 1, 0.3333333333333333, 0.003041, 0.009123, 0.0005003221135518961, 0.000500322113551896, 0.015176280796, 0.136586527164, 0.0002934466177152842, 0.0002934466177152842),
  @("Inflammatory-Mac", "Wnt11", "Fzd8", "ECs", 3, 1, 2.166324666666667, 6.498974, 0.2545974333859153, 0.2545974333859154, 3, 1, 2.157506, 6.472517999999999, 0.3549648016839517, 0.3549648016839516, 4.673858466281333, 42.064726196532, 0.09037312745107451, 0.09037312745107454),
  @("Inflammatory-Mac", "Wnt11", "Fzd8", "FAPs", 3, 1, 2.166324666666667, 6.498974, 0.2545974333859153, 0.2545974333859154, 3, 1, 2.913654666666667, 8.740964, 0.4793705560628122, 0.4793705560628121, 6.311921974548445, 56.807297770936, 0.122046513214371, 0.122046513214371),
  @("Inflammatory-Mac", "Wnt11", "Fzd8", "Inflammatory-Mac", 3, 1, 2.166324666666667, 6.498974, 0.2545974333859153, 0.2545974333859154, 2, 0.6666666666666666, 0.018986, 0.05695799999999999, 0.003123681567871193, 0.003123681567871192, 0.04112984012133333, 0.370168561092, 0.0007952813098948975, 0.0007952813098948977),
  @("Inflammatory-Mac", "Wnt11", "Fzd8", "MuSCs", 3, 1, 2.166324666666667, 6.498974, 0.2545974333859153, 0.2545974333859154, 3, 1, 0.9848966666666668, 2.95469, 0.1620406385718132, 0.1620406385718132, 2.133605943117778, 19.20245348806, 0.04125513068459839, 0.0412551306845984),
  @("Inflammatory-Mac", "Wnt11", "Fzd8", "Neutrophils", 3, 1, 2.166324666666667, 6.498974, 0.2545974333859153, 0.2545974333859154, 1, 0.3333333333333333, 0.003041, 0.009123, 0.0005003221135518961, 0.000500322113551896, 0.006587793311333333, 0.059290139802, 0.0001273807259765292, 0.0001273807259765292),
  @("MuSCs", "Wnt11", "Fzd8", "ECs", 3, 1, 0.318005, 0.9540150000000001, 0.03737355625852081, 0.03737355625852082, 3, 1, 2.157506, 6.472517999999999, 0.3549648016839517, 0.3549648016839516, 0.68609769553, 6.17487925977, 0.01326629698552985, 0.01326629698552985),
  @("MuSCs", "Wnt11", "Fzd8", "FAPs", 3, 1, 0.318005, 0.9540150000000001, 0.03737355625852081, 0.03737355625852082, 3, 1, 2.913654666666667, 8.740964, 0.4793705560628122, 0.4793705560628121, 0.9265567522733333, 8.33901077046, 0.01791578244569192, 0.01791578244569192),
  @("MuSCs", "Wnt11", "Fzd8", "Inflammatory-Mac", 3, 1, 0.318005, 0.9540150000000001, 0.03737355625852081, 0.03737355625852082, 2, 0.6666666666666666, 0.018986, 0.05695799999999999, 0.003123681567871193, 0.003123681567871192, 0.006037642930000001, 0.05433878637, 0.0001167430888105385, 0.0001167430888105385),
  @("MuSCs", "Wnt11", "Fzd8", "MuSCs", 3, 1, 0.318005, 0.9540150000000001, 0.03737355625852081, 0.03737355625852082, 3, 1, 0.9848966666666668, 2.95469, 0.1620406385718132, 0.1620406385718132, 0.3132020644833334, 2.81881858035, 0.006056034921830299, 0.006056034921830299),
  @("MuSCs", "Wnt11", "Fzd8", "Neutrophils", 3, 1, 0.318005, 0.9540150000000001, 0.03737355625852081, 0.03737355625852082, 1, 0.3333333333333333, 0.003041, 0.009123, 0.0005003221135518961, 0.000500322113551896, 0.0009670532050000001, 0.008703478845, 0.00001869881665821382, 0.00001869881665821382),
  @("Neutrophils", "Wnt11", "Fzd8", "ECs", 2, 0.6666666666666666, 0.2193053333333333, 0.657916, 0.0257738721502083, 0.02577387215020831, 3, 1, 2.157506, 6.472517999999999, 0.3549648016839517, 0.3549648016839516, 0.4731525724986666, 4.258373152487999, 0.009148817416426215, 0.009148817416426215),
  @("Neutrophils", "Wnt11", "Fzd8", "FAPs", 2, 0.6666666666666666, 0.2193053333333333, 0.657916, 0.0257738721502083, 0.02577387215020831, 3, 1, 2.913654666666667, 8.740964, 0.4793705560628122, 0.4793705560628121, 0.6389800078915555, 5.750820071023999, 0.01235523542453718, 0.01235523542453718),
  @("Neutrophils", "Wnt11", "Fzd8", "Inflammatory-Mac", 2, 0.6666666666666666, 0.2193053333333333, 0.657916, 0.0257738721502083, 0.02577387215020831, 2, 0.6666666666666666, 0.018986, 0.05695799999999999, 0.003123681567871193, 0.003123681567871192, 0.004163731058666667, 0.03747357952799999, 0.00008050936936827435, 0.00008050936936827435),
  @("Neutrophils", "Wnt11", "Fzd8", "MuSCs", 2, 0.6666666666666666, 0.2193053333333333, 0.657916, 0.0257738721502083, 0.02577387215020831, 3, 1, 0.9848966666666668, 2.95469, 0.1620406385718132, 0.1620406385718132, 0.2159930917822222, 1.94393782604, 0.004176414701688026, 0.004176414701688026),
  @("Neutrophils", "Wnt11", "Fzd8", "Neutrophils", 2, 0.6666666666666666, 0.2193053333333333, 0.657916, 0.0257738721502083, 0.02577387215020831, 1, 0.3333333333333333, 0.003041, 0.009123, 0.0005003221135518961, 0.000500322113551896, 0.0006669075186666666, 0.006002167667999999, 0.00001289523818860857, 0.00001289523818860857),
  @("Resolving-Mac", "Wnt11", "Fzd8", "ECs", 3, 1, 0.7241559999999999, 2.172468, 0.08510647633196142, 0.08510647633196144, 3, 1, 2.157506, 6.472517999999999, 0.3549648016839517, 0.3549648016839516, 1.562370914936, 14.061338234424, 0.03020980349319461, 0.03020980349319462),
  @("Resolving-Mac", "Wnt11", "Fzd8", "FAPs", 3, 1, 0.7241559999999999, 2.172468, 0.08510647633196142, 0.08510647633196144, 3, 1, 2.913654666666667, 8.740964, 0.4793705560628122, 0.4793705560628121, 2.109940508794666, 18.989464579152, 0.04079753888379891, 0.04079753888379892),
  @("Resolving-Mac", "Wnt11", "Fzd8", "Inflammatory-Mac", 3, 1, 0.7241559999999999, 2.172468, 0.08510647633196142, 0.08510647633196144, 2, 0.6666666666666666, 0.018986, 0.05695799999999999, 0.003123681567871193, 0.003123681567871192, 0.013748825816, 0.123739432344, 0.0002658455314246139, 0.0002658455314246139),
  @("Resolving-Mac", "Wnt11", "Fzd8", "MuSCs", 3, 1, 0.7241559999999999, 2.172468, 0.08510647633196142, 0.08510647633196144, 3, 1, 0.9848966666666668, 2.95469, 0.1620406385718132, 0.1620406385718132, 0.7132188305466667, 6.41896947492, 0.01379070777142794, 0.01379070777142794),
  @("Resolving-Mac", "Wnt11", "Fzd8", "Resolving-Mac", 3, 1, 0.7241559999999999, 2.172468, 0.08510647633196142, 0.08510647633196144, 1, 0.3333333333333333, 0.003041, 0.009123, 0.0005003221135518961, 0.000500322113551896, 0.002202158396, 0.019819425564, 0.00004258065211536135, 0.00004258065211536136)
)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($i = 0; $i -lt $data.Length; $i++) {
  $rowVals = $data[$i]
  $targetRow = $i + 2
  for ($j = 0; $j -lt $rowVals.Length; $j++) {
    $ws.Cells.Item($targetRow, $j + 1).Value = $rowVals[$j]
  }
}

Write-Output "Updated $($data.Length) data rows (2..$($data.Length + 1)) on $($ws.Name)."
